# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record was inserted as a new data row (row 37 on the
# worksheet, i.e. the 36th data row after the header), pushing all the
# subsequent rows down by one. The last existing row (old row 144) is
# therefore now row 145, and the sheet's dimension grows from R144 to R145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37; this shifts rows 37:144 down to 38:145 and
# carries over the existing number-format styling (e.g. the date style on
# column D) from the row that used to be there.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A37").Value = 8
$ws.Range("B37").Value = "Terminal La Palmera de La Serena"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 44526
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 100112031
$ws.Range("G37").Value = "Poroto verde"
$ws.Range("H37").Value = "Magnum"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 520
$ws.Range("K37").Value = 22000
$ws.Range("L37").Value = 23000
$ws.Range("M37").Value = 22500
$ws.Range("N37").Value = "$/malla 25 kilos"
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 900
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
